{"js": "// Replace the date line and each two-digit-multiplication answer cell with\n// its new value. Every old value in this document is unique, so we can\n// safely locate each one with a body-wide search-and-replace.\nconst replacements = [\n  [\"2025-06-07 Saturday\", \"2025-06-08 Sunday\"],\n  [\"71\u00d747=3337\", \"82\u00d748=3936\"],\n  [\"97\u00d716=1552\", \"26\u00d796=2496\"],\n  [\"79\u00d716=1264\", \"26\u00d745=1170\"],\n  [\"42\u00d778=3276\", \"83\u00d753=4399\"],\n  [\"94\u00d736=3384\", \"43\u00d780=3440\"],\n  [\"29\u00d763=1827\", \"32\u00d717=544\"],\n  [\"88\u00d798=8624\", \"55\u00d760=3300\"],\n  [\"19\u00d745=855\", \"99\u00d722=2178\"],\n  [\"87\u00d774=6438\", \"99\u00d760=5940\"],\n  [\"67\u00d757=3819\", \"33\u00d750=1650\"],\n  [\"88\u00d758=5104\", \"95\u00d746=4370\"],\n  [\"62\u00d785=5270\", \"30\u00d725=750\"],\n  [\"95\u00d772=6840\", \"66\u00d788=5808\"],\n  [\"92\u00d725=2300\", \"66\u00d733=2178\"],\n  [\"88\u00d752=4576\", \"88\u00d750=4400\"],\n  [\"52\u00d777=4004\", \"14\u00d798=1372\"],\n  [\"64\u00d740=2560\", \"83\u00d739=3237\"],\n  [\"44\u00d729=1276\", \"49\u00d720=980\"],\n  [\"87\u00d748=4176\", \"28\u00d763=1764\"],\n  [\"70\u00d782=5740\", \"23\u00d760=1380\"],\n  [\"64\u00d722=1408\", \"31\u00d774=2294\"],\n  [\"51\u00d768=3468\", \"97\u00d779=7663\"],\n  [\"56\u00d757=3192\", \"16\u00d785=1360\"],\n  [\"19\u00d734=646\", \"51\u00d715=765\"],\n  [\"89\u00d769=6141\", \"62\u00d794=5828\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and each two-digit-multiplication answer cell with\n# its new value. Every old value in this document is unique, so a simple\n# Find/Replace over the whole document body is safe for each pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-06-07 Saturday\", \"2025-06-08 Sunday\"),\n    @(\"71\u00d747=3337\", \"82\u00d748=3936\"),\n    @(\"97\u00d716=1552\", \"26\u00d796=2496\"),\n    @(\"79\u00d716=1264\", \"26\u00d745=1170\"),\n    @(\"42\u00d778=3276\", \"83\u00d753=4399\"),\n    @(\"94\u00d736=3384\", \"43\u00d780=3440\"),\n    @(\"29\u00d763=1827\", \"32\u00d717=544\"),\n    @(\"88\u00d798=8624\", \"55\u00d760=3300\"),\n    @(\"19\u00d745=855\", \"99\u00d722=2178\"),\n    @(\"87\u00d774=6438\", \"99\u00d760=5940\"),\n    @(\"67\u00d757=3819\", \"33\u00d750=1650\"),\n    @(\"88\u00d758=5104\", \"95\u00d746=4370\"),\n    @(\"62\u00d785=5270\", \"30\u00d725=750\"),\n    @(\"95\u00d772=6840\", \"66\u00d788=5808\"),\n    @(\"92\u00d725=2300\", \"66\u00d733=2178\"),\n    @(\"88\u00d752=4576\", \"88\u00d750=4400\"),\n    @(\"52\u00d777=4004\", \"14\u00d798=1372\"),\n    @(\"64\u00d740=2560\", \"83\u00d739=3237\"),\n    @(\"44\u00d729=1276\", \"49\u00d720=980\"),\n    @(\"87\u00d748=4176\", \"28\u00d763=1764\"),\n    @(\"70\u00d782=5740\", \"23\u00d760=1380\"),\n    @(\"64\u00d722=1408\", \"31\u00d774=2294\"),\n    @(\"51\u00d768=3468\", \"97\u00d779=7663\"),\n    @(\"56\u00d757=3192\", \"16\u00d785=1360\"),\n    @(\"19\u00d734=646\", \"51\u00d715=765\"),\n    @(\"89\u00d769=6141\", \"62\u00d794=5828\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n$d.Save()\n"}
